# "update scripts wuth new tpm"
# Recomputed cell-cell communication (NATMI) scores for the Spon2-Itgam
# ligand-receptor pair using the updated TPM values. This refreshes the
# numeric score columns (E:T) for the existing sending/target cluster
# combinations (rows 2-7) and adds the new "Resolving-Mac" cluster as
# both a sending cluster and a target cluster (rows 8-9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Spon2"
$ws.Range("C2").Value = "Itgam"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4500866666666667
$ws.Range("H2").Value = 1.35026
$ws.Range("I2").Value = 0.02628438542510526
$ws.Range("J2").Value = 0.02628438542510525
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02725333333333333
$ws.Range("N2").Value = 0.08176
$ws.Range("O2").Value = 0.0007089206372884383
$ws.Range("P2").Value = 0.0007089206372884382
$ws.Range("Q2").Value = 0.01226636195555556
$ws.Range("R2").Value = 0.1103972576
$ws.Range("S2").Value = [double]"1.863354326630056e-05"
$ws.Range("T2").Value = [double]"1.863354326630055e-05"

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Spon2"
$ws.Range("C3").Value = "Itgam"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4500866666666667
$ws.Range("H3").Value = 1.35026
$ws.Range("I3").Value = 0.02628438542510526
$ws.Range("J3").Value = 0.02628438542510525
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 38.416166
$ws.Range("N3").Value = 115.248498
$ws.Range("O3").Value = 0.9992910793627116
$ws.Range("P3").Value = 0.9992910793627116
$ws.Range("Q3").Value = 17.29060410105333
$ws.Range("R3").Value = 155.61543690948
$ws.Range("S3").Value = 0.02626575188183896
$ws.Range("T3").Value = 0.02626575188183895

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Spon2"
$ws.Range("C4").Value = "Itgam"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 15.76143266666667
$ws.Range("H4").Value = 47.284298
$ws.Range("I4").Value = 0.9204439983318276
$ws.Range("J4").Value = 0.9204439983318274
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02725333333333333
$ws.Range("N4").Value = 0.08176
$ws.Range("O4").Value = 0.0007089206372884383
$ws.Range("P4").Value = 0.0007089206372884382
$ws.Range("Q4").Value = 0.4295515782755556
$ws.Range("R4").Value = 3.86596420448
$ws.Range("S4").Value = 0.0006525217458857174
$ws.Range("T4").Value = 0.0006525217458857172

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Spon2"
$ws.Range("C5").Value = "Itgam"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 15.76143266666667
$ws.Range("H5").Value = 47.284298
$ws.Range("I5").Value = 0.9204439983318276
$ws.Range("J5").Value = 0.9204439983318274
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 38.416166
$ws.Range("N5").Value = 115.248498
$ws.Range("O5").Value = 0.9992910793627116
$ws.Range("P5").Value = 0.9992910793627116
$ws.Range("Q5").Value = 605.4938137204892
$ws.Range("R5").Value = 5449.444323484404
$ws.Range("S5").Value = 0.9197914765859418
$ws.Range("T5").Value = 0.9197914765859417

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Spon2"
$ws.Range("C6").Value = "Itgam"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9026056666666666
$ws.Range("H6").Value = 2.707817
$ws.Range("I6").Value = 0.05271081546417152
$ws.Range("J6").Value = 0.05271081546417151
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.02725333333333333
$ws.Range("N6").Value = 0.08176
$ws.Range("O6").Value = 0.0007089206372884383
$ws.Range("P6").Value = 0.0007089206372884382
$ws.Range("Q6").Value = 0.02459901310222222
$ws.Range("R6").Value = 0.22139111792
$ws.Range("S6").Value = [double]"3.736778489085374e-05"
$ws.Range("T6").Value = [double]"3.736778489085373e-05"

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Spon2"
$ws.Range("C7").Value = "Itgam"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.9026056666666666
$ws.Range("H7").Value = 2.707817
$ws.Range("I7").Value = 0.05271081546417152
$ws.Range("J7").Value = 0.05271081546417151
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 38.416166
$ws.Range("N7").Value = 115.248498
$ws.Range("O7").Value = 0.9992910793627116
$ws.Range("P7").Value = 0.9992910793627116
$ws.Range("Q7").Value = 34.67464912320733
$ws.Range("R7").Value = 312.071842108866
$ws.Range("S7").Value = 0.05267344767928066
$ws.Range("T7").Value = 0.05267344767928065

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Spon2"
$ws.Range("C8").Value = "Itgam"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.009603
$ws.Range("H8").Value = 0.028809
$ws.Range("I8").Value = 0.0005608007788958107
$ws.Range("J8").Value = 0.0005608007788958106
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02725333333333333
$ws.Range("N8").Value = 0.08176
$ws.Range("O8").Value = 0.0007089206372884383
$ws.Range("P8").Value = 0.0007089206372884382
$ws.Range("Q8").Value = 0.00026171376
$ws.Range("R8").Value = 0.00235542384
$ws.Range("S8").Value = [double]"3.975632455666707e-07"
$ws.Range("T8").Value = [double]"3.975632455666705e-07"

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Spon2"
$ws.Range("C9").Value = "Itgam"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.009603
$ws.Range("H9").Value = 0.028809
$ws.Range("I9").Value = 0.0005608007788958107
$ws.Range("J9").Value = 0.0005608007788958106
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 38.416166
$ws.Range("N9").Value = 115.248498
$ws.Range("O9").Value = 0.9992910793627116
$ws.Range("P9").Value = 0.9992910793627116
$ws.Range("Q9").Value = 0.368910442098
$ws.Range("R9").Value = 3.320193978882
$ws.Range("S9").Value = 0.000560403215650244
$ws.Range("T9").Value = 0.0005604032156502439

